$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting existing rows 20:83 down to 21:84
# (preserves formatting of the row being pushed down, including the date style on column D)
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the latest weekly price entry
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44481
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100108
$ws.Range("H20").Value = "Tropicales y subtropicales"
$ws.Range("I20").Value = 100108002
$ws.Range("J20").Value = "Mango"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 7500
$ws.Range("O20").Value = 8000
$ws.Range("P20").Value = 7750
$ws.Range("Q20").Value = "$/bandeja 4 kilos"
$ws.Range("R20").Value = "Perú"
$ws.Range("S20").Value = 1938
$ws.Range("T20").Value = 4
